$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two accuracy values in row 5 (GLCM balanced-without-background row)
$ws.Range("B5").Value = 0.28057553956834502
$ws.Range("C5").Value = 0.36690647482014299

# Row 2's explicit row height (wrapped header) is no longer needed -- clear it
# back to the sheet default by auto-fitting the row.
$ws.Rows("2").AutoFit()

# Configure the print setup for the sheet (portrait orientation).
$ws.PageSetup.Orientation = 1

# Move/save the current selection to B10, matching the last cursor position.
$ws.Range("B10").Select()
